$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5258
$ws1.Range("F3").Value = 377
$ws1.Range("F6").Value = 797
$ws1.Range("F7").Value = 294

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 8

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5258
$ws4.Range("F3").Value = 377
$ws4.Range("F6").Value = 797
$ws4.Range("F8").Value = 294
$ws4.Range("F10").Value = 8
